# Update countries & provincias Spain
#
# Applies the daily COVID-19 data refresh to the "Pais" sheet:
#  - straightforward numeric updates for a handful of existing country rows
#  - "Kenia" re-ranked ahead of "Jamaica" (its row physically relocates from
#    122 to 119, pushing Jamaica/Taiwan/Reunion down one row each) together
#    with Kenia's own stats refresh
#  - "Namibia" swaps places with "San Vicente y las Granadinas" (rows 193/194)
#  - "Comoras" swaps places with "San Pedro y Miquelon" (rows 217/218)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain per-cell numeric refreshes (country stays in place) ---------
$numericUpdates = @{
    4   = @{ B = 1134084; C = 3054;  E = 906414; G = 135; H = 65888 }
    9   = @{ B = 164271;  C = 194;   E = 28535 }
    27  = @{ B = 18770;   C = 678;   D = 4753;  E = 13585; G = 15; H = 432 }
    57  = @{ D = 1320;    E = 2983;  G = 4;     H = 229 }
    77  = @{ D = 1706;    E = 82 }
    79  = @{ B = 1594;    C = 39;    E = 1235;  G = 4; H = 72 }
    99  = @{ D = 197;     E = 511 }
    167 = @{ D = 32;      E = 28 }
}

foreach ($row in $numericUpdates.Keys) {
    $rowValues = $numericUpdates[$row]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}

# --- Kenia moves ahead of Jamaica (row 122 -> row 119) ------------------
# Row 119 becomes Kenia (refreshed stats); Jamaica/Taiwan/Reunion each shift
# down one row, keeping their own (unchanged) statistics.
$rows119to122 = @{
    119 = @{ Country = "Kenia";   B = 435; C = 24; D = 152; E = 261; F = 2; G = 1; H = 22 }
    120 = @{ Country = "Jamaica"; B = 432; C = 10; D = 31;  E = 393; F = 1; G = 0; H = 8 }
    121 = @{ Country = "Taiwan";  B = 432; C = 3;  D = 324; E = 102; F = 0; G = 0; H = 6 }
    122 = @{ Country = "Reunion"; B = 422; C = 0;  D = 300; E = 122; F = 2; G = 0; H = 0 }
}

foreach ($row in $rows119to122.Keys) {
    $rowValues = $rows119to122[$row]
    $ws.Range("A$row").Value = $rowValues["Country"]
    foreach ($col in @("B", "C", "D", "E", "F", "G", "H")) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}

# --- Namibia <-> San Vicente y las Granadinas (rows 193 / 194) ---------
$ws.Range("A193").Value = "Namibia"
$ws.Range("A194").Value = "San Vicente y las Granadinas"

# --- Comoras <-> San Pedro y Miquelon (rows 217 / 218) ------------------
$ws.Range("A217").Value = "Comoras"
$ws.Range("A218").Value = "San Pedro y Miquelon"
